$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "status" column (G) with "Done!" for the newly collected experiment
# results on rows 7-10, matching the same value already used in rows 3-6.
$ws.Range("G7").Value = "Done!"
$ws.Range("G8").Value = "Done!"
$ws.Range("G9").Value = "Done!"
$ws.Range("G10").Value = "Done!"

# Update the selected cell to reflect where the user finished editing.
$ws.Range("A9").Select()
